$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1, J1 (same style as other header cells, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I and J (values mirror each other per row)
$values = @{
    2 = 9
    3 = 8
    4 = 9
    5 = 9
    6 = 9
    7 = 8
    8 = 8
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 9).Value = $values[$row]
    $ws.Cells.Item($row, 10).Value = $values[$row]
}
